$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "Ativação:" date value (row 8, columns B and C): 01/01/2012 -> 01/01/2021 ---
# Assigning a date-shaped literal straight to .Value gets auto-coerced into a
# real date serial by Excel's type inference, which would also bolt a new
# number format onto the cell's style. Routing the text through a
# text-producing formula and then "paste values" keeps it a plain shared
# string using the cell's existing style (no style/numFmt churn).
$ws.Range("B8").Formula = '="01/01/2021"'
$ws.Range("C8").Formula = '="01/01/2021"'
$wb.Application.Calculate()
$ws.Range("B8:C8").Copy()
$ws.Range("B8:C8").PasteSpecial(-4163) # xlPasteValues

# --- 2. "Docentes responsáveis:" value (row 13, columns B and C) ---
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

# --- 3. "Critério:" value (row 20, columns B and C) ---
$ws.Range("B20").Value = "Média Aritmética das atividades avaliativas realizadas."
$ws.Range("C20").Value = "Média Aritmética das atividades avaliativas realizadas."

# --- 4. "Norma de recuperação:" value, trailing period removed (row 21, B and C) ---
$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação"
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação"

# --- 5. New requisite row 25, mirroring row 24's layout/style ---
$ws.Range("B24:C24").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B25").Value = "LOQ4240 -  Administração e Organização II  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOQ4240 -  Administração e Organização II  (Requisito fraco)`n"
$ws.Rows(25).RowHeight = 30

$wb.Save()
